# Intermediate check-in of graph prompt types
# Adds a new "graphExample" test-form entry to the "survey" sheet (mirrors the
# existing pattern used for every other form in that section) and a matching
# row in the "choices" sheet's test_forms list.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "survey" sheet: insert the new graphExample block right after the last
#    existing form block (currently rows 47-49, "Adult Coverage Immunizations").
#    The new block follows the exact same 3-row shape:
#      row 50 -> form name (col A)
#      row 51 -> hash-string formula (col B), external_link/Open form (E/G)
#      row 52 -> "exit section" marker (col C)
# ---------------------------------------------------------------------------
$wsSurvey = $wb.Worksheets.Item("survey")

# Carry over the formatting from the previous block so the new rows pick up
# identical cell styles (plain style for A/B "marker" cells, quote-prefix +
# wrap style for the formula cell) instead of Excel synthesizing new ones.
$wsSurvey.Range("A48").Copy()
$wsSurvey.Range("A51").PasteSpecial(-4122)   # xlPasteFormats
$wsSurvey.Range("B48").Copy()
$wsSurvey.Range("B51").PasteSpecial(-4122)   # xlPasteFormats
$wsSurvey.Range("A49").Copy()
$wsSurvey.Range("A52").PasteSpecial(-4122)   # xlPasteFormats
$wsSurvey.Range("B49").Copy()
$wsSurvey.Range("B52").PasteSpecial(-4122)   # xlPasteFormats

$wsSurvey.Range("A50").Value = "graphExample"
# Doubled leading apostrophe so the literal value keeps a single leading "'"
# (Excel's text-prefix escape) instead of that char being swallowed.
$wsSurvey.Range("B51").Value = "''?' + opendatakit.getHashString('../tables/graphExample/forms/graphExample/',null)"
$wsSurvey.Range("E51").Value = "external_link"
$wsSurvey.Range("G51").Value = "Open form"
$wsSurvey.Range("C52").Value = "exit section"

$wsSurvey.Rows.Item(50).RowHeight = 17.5
$wsSurvey.Rows.Item(51).RowHeight = 66
$wsSurvey.Rows.Item(52).RowHeight = 17

# ---------------------------------------------------------------------------
# 2. "choices" sheet: append the matching test_forms choice row.
# ---------------------------------------------------------------------------
$wsChoices = $wb.Worksheets.Item("choices")

$wsChoices.Range("A15").Copy()
$wsChoices.Range("A16").PasteSpecial(-4122)  # xlPasteFormats

$wsChoices.Range("A16").Value = "test_forms"
$wsChoices.Range("B16").Value = "graphExample"
$wsChoices.Range("C16").Value = "Graph Examples"

# ---------------------------------------------------------------------------
# 3. Restore/update the on-screen selection state for both touched sheets.
#    "choices" stays the active tab (as in the original file); "survey" keeps
#    a (non-active) selection pointing at the newly-added formula cell.
# ---------------------------------------------------------------------------
$wsSurvey.Activate()
$wsSurvey.Range("B51").Select()

$wsChoices.Activate()
$wsChoices.Range("C23").Select()
